$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.299.94'
$ws.Range('E2').Value = '  -2.16%  '

$ws.Range('D3').Value = '3.371.71'
$ws.Range('E3').Value = '  -2.32%  '

$ws.Range('E4').Value = '  +0.01%  '

$c = $ws.Range('D5')
$origStyle = $c.Style
$c.Value = "'566.40"
$c.Style = $origStyle
$ws.Range('E5').Value = '  -2.41%  '

$c = $ws.Range('D6')
$origStyle = $c.Style
$c.Value = "'140.19"
$c.Style = $origStyle
$ws.Range('E6').Value = '  -6.74%  '

$ws.Range('D8').Value = '3.371.66'
$ws.Range('E8').Value = '  -2.36%  '

$ws.Range('E9').Value = '  -0.73%  '

$c = $ws.Range('D10')
$origStyle = $c.Style
$c.Value = "'7.48"
$c.Style = $origStyle
$ws.Range('E10').Value = '  -4.38%  '

$ws.Range('E11').Value = '  -2.55%  '

$c = $ws.Range('D12')
$origStyle = $c.Style
$c.Value = "'0.386"
$c.Style = $origStyle
$ws.Range('E12').Value = '  -1.47%  '

$ws.Range('D13').Value = '3.948.12'
$ws.Range('E13').Value = '  -2.30%  '

$ws.Range('E14').Value = '  +1.07%  '

$c = $ws.Range('D15')
$origStyle = $c.Style
$c.Value = "'27.92"
$c.Style = $origStyle
$ws.Range('E15').Value = '  -0.24%  '

$ws.Range('D16').Value = '3.395.13'
$ws.Range('E16').Value = '  -1.49%  '

$ws.Range('E17').Value = '  -3.89%  '

$ws.Range('D18').Value = '60.399.39'
$ws.Range('E18').Value = '  -2.18%  '

$c = $ws.Range('D19')
$origStyle = $c.Style
$c.Value = "'6.16"
$c.Style = $origStyle
$ws.Range('E19').Value = '  -2.38%  '

$c = $ws.Range('D20')
$origStyle = $c.Style
$c.Value = "'13.80"
$c.Style = $origStyle
$ws.Range('E20').Value = '  -4.20%  '

$c = $ws.Range('D21')
$origStyle = $c.Style
$c.Value = "'8.98"
$c.Style = $origStyle
$ws.Range('E21').Value = '  -5.69%  '

$c = $ws.Range('D22')
$origStyle = $c.Style
$c.Value = "'384.75"
$c.Style = $origStyle
$ws.Range('E22').Value = '  -1.35%  '

$c = $ws.Range('D23')
$origStyle = $c.Style
$c.Value = "'0.554"
$c.Style = $origStyle
$ws.Range('E23').Value = '  -2.54%  '

$c = $ws.Range('D24')
$origStyle = $c.Style
$c.Value = "'72.97"
$c.Style = $origStyle
$ws.Range('E24').Value = '  +0.03%  '

$ws.Range('E25').Value = '  +0.08%  '

$ws.Range('E26').Value = '  -7.94%  '

$ws.Range('D27').Value = '3.525.39'
$ws.Range('E27').Value = '  -1.71%  '

$ws.Range('E28').Value = '  -2.38%  '

$c = $ws.Range('D29')
$origStyle = $c.Style
$c.Value = "'0.998"
$c.Style = $origStyle
$ws.Range('E29').Value = '  -0.16%  '

$c = $ws.Range('D30')
$origStyle = $c.Style
$c.Value = "'7.33"
$c.Style = $origStyle
$ws.Range('E30').Value = '  -5.67%  '

$c = $ws.Range('D31')
$origStyle = $c.Style
$c.Value = "'7.90"
$c.Style = $origStyle
$ws.Range('E31').Value = '  -4.52%  '

$ws.Range('E32').Value = '  -2.34%  '

$c = $ws.Range('D33')
$origStyle = $c.Style
$c.Value = "'1.40"
$c.Style = $origStyle
$ws.Range('E33').Value = '  -9.20%  '

$ws.Range('E34').Value = '  -0.03%  '

$c = $ws.Range('D35')
$origStyle = $c.Style
$c.Value = "'23.46"
$c.Style = $origStyle
$ws.Range('E35').Value = '  -2.45%  '

$ws.Range('D36').Value = '3.402.83'
$ws.Range('E36').Value = '  -2.09%  '

$c = $ws.Range('D37')
$origStyle = $c.Style
$c.Value = "'6.87"
$c.Style = $origStyle
$ws.Range('E37').Value = '  -2.70%  '

$c = $ws.Range('D38')
$origStyle = $c.Style
$c.Value = "'167.86"
$c.Style = $origStyle
$ws.Range('E38').Value = '  +0.61%  '

$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D39')
$origStyle = $c.Style
$c.Value = "'4.90"
$c.Style = $origStyle
$ws.Range('E39').Value = '  -6.82%  '

$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D40')
$origStyle = $c.Style
$c.Value = "'1.49"
$c.Style = $origStyle
$ws.Range('E40').Value = '  -5.08%  '

$c = $ws.Range('D41')
$origStyle = $c.Style
$c.Value = "'0.0767"
$c.Style = $origStyle
$ws.Range('E41').Value = '  -2.96%  '

$c = $ws.Range('D42')
$origStyle = $c.Style
$c.Value = "'26.96"
$c.Style = $origStyle
$ws.Range('E42').Value = '  -0.72%  '

$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range('D43')
$origStyle = $c.Style
$c.Value = "'1.00"
$c.Style = $origStyle
$ws.Range('E43').Value = '  +0.04%  '

$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range('D44')
$origStyle = $c.Style
$c.Value = "'0.778"
$c.Style = $origStyle
$ws.Range('E44').Value = '  -2.02%  '

$c = $ws.Range('D45')
$origStyle = $c.Style
$c.Value = "'4.41"
$c.Style = $origStyle
$ws.Range('E45').Value = '  -2.29%  '

$ws.Range('E46').Value = '  -2.64%  '

$ws.Range('E47').Value = '  -2.41%  '

$ws.Range('D48').Value = '2.518.98'
$ws.Range('E48').Value = '  -3.27%  '

$ws.Range('E49').Value = '  -5.27%  '

$c = $ws.Range('D50')
$origStyle = $c.Style
$c.Value = "'23.12"
$c.Style = $origStyle
$ws.Range('E50').Value = '  -0.49%  '

$c = $ws.Range('D51')
$origStyle = $c.Style
$c.Value = "'6.75"
$c.Style = $origStyle
$ws.Range('E51').Value = '  -3.28%  '
